$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.08397466666667
$ws.Range("H2").Value = 63.251924
$ws.Range("I2").Value = 0.06331801375981215
$ws.Range("J2").Value = 0.06331801375981214
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 58.18216281554355
$ws.Range("R2").Value = 523.6394653398919
$ws.Range("S2").Value = 0.01629848310804768
$ws.Range("T2").Value = 0.01629848310804768

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.08397466666667
$ws.Range("H3").Value = 63.251924
$ws.Range("I3").Value = 0.06331801375981215
$ws.Range("J3").Value = 0.06331801375981214
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 152.7581825222493
$ws.Range("R3").Value = 1374.823642700244
$ws.Range("S3").Value = 0.04279192345166321
$ws.Range("T3").Value = 0.0427919234516632

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 21.08397466666667
$ws.Range("H4").Value = 63.251924
$ws.Range("I4").Value = 0.06331801375981215
$ws.Range("J4").Value = 0.06331801375981214
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 10.96432745787289
$ws.Range("R4").Value = 98.67894712085599
$ws.Range("S4").Value = 0.003071420813794566
$ws.Range("T4").Value = 0.003071420813794564

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 21.08397466666667
$ws.Range("H5").Value = 63.251924
$ws.Range("I5").Value = 0.06331801375981215
$ws.Range("J5").Value = 0.06331801375981214
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 4.127342656814222
$ws.Range("R5").Value = 37.146083911328
$ws.Range("S5").Value = 0.001156186386306698
$ws.Range("T5").Value = 0.001156186386306698

$ws.Range("I6").Value = 0.8174956765497907
$ws.Range("J6").Value = 0.8174956765497907
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 751.1869644940037
$ws.Range("R6").Value = 6760.682680446035
$ws.Range("S6").Value = 0.2104288919373125
$ws.Range("T6").Value = 0.2104288919373125

$ws.Range("I7").Value = 0.8174956765497907
$ws.Range("J7").Value = 0.8174956765497907
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.5524843616491872
$ws.Range("T7").Value = 0.5524843616491872

$ws.Range("I8").Value = 0.8174956765497907
$ws.Range("J8").Value = 0.8174956765497907
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 141.5598778427922
$ws.Range("R8").Value = 1274.03890058513
$ws.Range("S8").Value = 0.03965495894528114
$ws.Range("T8").Value = 0.03965495894528114

$ws.Range("I9").Value = 0.8174956765497907
$ws.Range("J9").Value = 0.8174956765497907
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 53.28791251071556
$ws.Range("R9").Value = 479.5912125964401
$ws.Range("S9").Value = 0.0149274640180099
$ws.Range("T9").Value = 0.0149274640180099

$ws.Range("G10").Value = 39.60693866666666
$ws.Range("H10").Value = 118.820816
$ws.Range("I10").Value = 0.1189449677837485
$ws.Range("J10").Value = 0.1189449677837485
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 109.2971031582809
$ws.Range("R10").Value = 983.6739284245278
$ws.Range("S10").Value = 0.03061723565057786
$ws.Range("T10").Value = 0.03061723565057785

$ws.Range("G11").Value = 39.60693866666666
$ws.Range("H11").Value = 118.820816
$ws.Range("I11").Value = 0.1189449677837485
$ws.Range("J11").Value = 0.1189449677837485
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 286.9612614150773
$ws.Range("R11").Value = 2582.651352735696
$ws.Range("S11").Value = 0.08038603320171191
$ws.Range("T11").Value = 0.08038603320171189

$ws.Range("G12").Value = 39.60693866666666
$ws.Range("H12").Value = 118.820816
$ws.Range("I12").Value = 0.1189449677837485
$ws.Range("J12").Value = 0.1189449677837485
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 20.59684912407822
$ws.Range("R12").Value = 185.371642116704
$ws.Range("S12").Value = 0.005769764843429178
$ws.Range("T12").Value = 0.005769764843429176

$ws.Range("G13").Value = 39.60693866666666
$ws.Range("H13").Value = 118.820816
$ws.Range("I13").Value = 0.1189449677837485
$ws.Range("J13").Value = 0.1189449677837485
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 7.753348694883555
$ws.Range("R13").Value = 69.780138253952
$ws.Range("S13").Value = 0.002171934088029529
$ws.Range("T13").Value = 0.002171934088029528

$ws.Range("G14").Value = 0.08036333333333333
$ws.Range("H14").Value = 0.24109
$ws.Range("I14").Value = 0.0002413419066486121
$ws.Range("J14").Value = 0.0002413419066486121
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.2217661811077777
$ws.Range("R14").Value = 1.99589562997
$ws.Range("S14").Value = 0.00006212303190206854
$ws.Range("T14").Value = 0.00006212303190206854

$ws.Range("G15").Value = 0.08036333333333333
$ws.Range("H15").Value = 0.24109
$ws.Range("I15").Value = 0.0002413419066486121
$ws.Range("J15").Value = 0.0002413419066486121
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 0.5822505924766667
$ws.Range("R15").Value = 5.24025533229
$ws.Range("S15").Value = 0.0001631049962205336
$ws.Range("T15").Value = 0.0001631049962205336

$ws.Range("G16").Value = 0.08036333333333333
$ws.Range("H16").Value = 0.24109
$ws.Range("I16").Value = 0.0002413419066486121
$ws.Range("J16").Value = 0.0002413419066486121
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.04179145138444445
$ws.Range("R16").Value = 0.37612306246
$ws.Range("S16").Value = 0.00001170697738771917
$ws.Range("T16").Value = 0.00001170697738771917

$ws.Range("G17").Value = 0.08036333333333333
$ws.Range("H17").Value = 0.24109
$ws.Range("I17").Value = 0.0002413419066486121
$ws.Range("J17").Value = 0.0002413419066486121
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.01573171183111111
$ws.Range("R17").Value = 0.14158540648
$ws.Range("S17").Value = 0.000004406901138290778
$ws.Range("T17").Value = 0.000004406901138290778
